# Dustin Bramos effort log - add Sprint #4 entry and update the Count summary.
$wb = $excel.ActiveWorkbook

$wsEvidence = $wb.Worksheets.Item("Evidence")
$wsCount = $wb.Worksheets.Item("Count")

# --- Evidence sheet: fill in the Sprint # 4 row (row 7) ---
$wsEvidence.Range("C7").Value = 35
$wsEvidence.Range("D7").Value = "1, 2, 3, 5, 6[CS]"
$wsEvidence.Range("E7").Value = "Created the History page's query and helped with controllers. Created a MySQL script to completely fill the database with test data. Made a small change to the ERD. Collaborated on the powerpoint presentation."
$wsEvidence.Range("F7").Value = "userQueries.php, userControllers.php, Sprint4.pptx, Test Data Insertion Script.sql and GPTMS_ERD.mwb"
$wsEvidence.Range("G7").Value = "userQueries.php, userControllers.php, Sprint4.pptx, Test Data Insertion Script.sql and GPTMS_ERD.mwb"
$wsEvidence.Range("H7").Value = "The query inside userQueries.php, the script inside the test insertion script, the Sprint 4 powerpoint, and the database ERD."

# The long wrapped description text now needs more vertical room, same as the
# other filled-in sprint rows above it.
$wsEvidence.Rows.Item(7).RowHeight = 93.6

# --- Count sheet: update the Sprint 4 objective tally row (row 8) ---
$wsCount.Range("C8").Value = 10
$wsCount.Range("D8").Value = 1
$wsCount.Range("E8").Value = 1
$wsCount.Range("G8").Value = 3
$wsCount.Range("H8").Value = 20

# --- Selection / active sheet bookkeeping to match the final saved state ---
$wsCount.Activate()
$wsCount.Range("C8").Select()

$wsEvidence.Activate()
$wsEvidence.Range("H7").Select()
